$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.897.76"
$ws.Range("E2").Value = "  -1.39%  "
$ws.Range("D3").Value = "3.513.75"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "608.87"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.92"
$ws.Range("E6").Value = "  -2.50%  "
$ws.Range("D7").Value = "3.514.31"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("E9").Value = "  -2.03%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("E11").Value = "  +4.87%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.422"
$ws.Range("E12").Value = "  -2.27%  "
$ws.Range("E13").Value = "  +0.78%  "
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "31.97"
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "4.108.09"
$ws.Range("E15").Value = "  +0.34%  "
$ws.Range("D16").Value = "3.507.84"
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "66.980.74"
$ws.Range("E17").Value = "  -1.18%  "
$ws.Range("E18").Value = "  -0.44%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.76"
$ws.Range("E19").Value = "  +8.14%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.46"
$ws.Range("E20").Value = "  -0.82%  "
$ws.Range("E21").Value = "  -0.51%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "438.10"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("E23").Value = "  -3.00%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.57"
$ws.Range("E24").Value = "  +0.18%  "
$ws.Range("D25").Value = "3.653.45"
$ws.Range("E25").Value = "  +0.33%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000121"
$ws.Range("E27").Value = "  -4.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.79"
$ws.Range("E28").Value = "  -1.84%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.29"
$ws.Range("E29").Value = "  -5.07%  "
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.62"
$ws.Range("E31").Value = "  -1.92%  "
$ws.Range("E32").Value = "  -2.06%  "
$ws.Range("E33").Value = "  -0.11%  "
$ws.Range("E34").Value = "  -0.48%  "
$ws.Range("E35").Value = "  -3.11%  "
$ws.Range("E36").Value = "  -2.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "8.05"
$ws.Range("E37").Value = "  +0.05%  "
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "175.44"
$ws.Range("E40").Value = "  -1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0895"
$ws.Range("E41").Value = "  -0.57%  "
$ws.Range("E42").Value = "  -0.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.07"
$ws.Range("E43").Value = "  -11.60%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.895"
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "46.17"
$ws.Range("E45").Value = "  -1.01%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "28.04"
$ws.Range("E46").Value = "  -8.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.26"
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.47"
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.992"
$ws.Range("E50").Value = "  -0.13%  "
$ws.Range("E51").Value = "  -1.61%  "
